function RGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

# Swap the deck's theme colour scheme from "Integral" (Red Violet) to the
# default "Office Theme" palette - i.e. applying the Office Theme design.
$cs.Colors(1).RGB  = RGB 0x00 0x00 0x00   # dk1
$cs.Colors(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1
$cs.Colors(3).RGB  = RGB 0x44 0x54 0x6A   # dk2
$cs.Colors(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2
$cs.Colors(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1
$cs.Colors(6).RGB  = RGB 0xED 0x7D 0x31   # accent2
$cs.Colors(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3
$cs.Colors(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4
$cs.Colors(9).RGB  = RGB 0x44 0x72 0xC4   # accent5
$cs.Colors(10).RGB = RGB 0x70 0xAD 0x47   # accent6
$cs.Colors(11).RGB = RGB 0x05 0x63 0xC1   # hlink
$cs.Colors(12).RGB = RGB 0x95 0x4F 0x72   # folHlink
